$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Runmode" column (C) flag flips: the Y flag moves from
# Web_FILE_UPLOAD_EXISTING_USER (row 8) to Web_MESSAGING (row 5).
$ws.Range("C5").Value = "Y"
$ws.Range("C8").Value = "N"

# Active selection moves from C7 to C6.
[void]$ws.Range("C6").Select()
